$d = $word.ActiveDocument

# The author corrected the misspelled/abbreviated first name "Mohd" to
# "Mohammad" in the "Secretary/Co-Leader" line. Word tracks the location
# of the most recent edit with the hidden "_GoBack" bookmark, so we place
# that bookmark where the edit happens (immediately after "Mohd") BEFORE
# performing the text replacement - this also matches how real Word keeps
# the untouched runs around it ("-" / " " / "Secretary/Co-Leader") intact
# instead of merging them together.
$findRange = $d.Content
$findRange.Find.Execute("Mohd", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

if ($findRange.Find.Found) {
    $editPoint = $findRange.End
    $goBackRange = $d.Range($editPoint, $editPoint)
    $d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
}

# Now perform the actual text correction: "Mohd" -> "Mohammad".
$replaceRange = $d.Content
$replaceRange.Find.Execute("Mohd", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "Mohammad", 2)
